# The post previously stored in row 579 ("大切な時間を大切じゃないことに使わないこと")
# was removed. Deleting the entire row shifts every subsequent row up by one,
# which also updates the sheet's used-range dimension automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(579).Delete()
